$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "K" column (column G) values per regenerated save_data
# (std/mean recalculated and s_vals rewritten upstream; here we just
# write the resulting K values into the sheet).
$kValues = @{
    2  = 0
    3  = 2
    4  = 0
    5  = 1
    6  = 0
    7  = 1
    8  = 0
    9  = 1
    10 = 0
    11 = 1
    12 = 2
    13 = 1
    14 = 0
    15 = 2
    16 = 1
    17 = 0
    18 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 0
    23 = 0
    24 = 2
    25 = 3
    26 = 0
    27 = 1
    28 = 0
    29 = 0
    30 = 1
    31 = 3
    32 = 0
    33 = 1
    34 = 1
    35 = 3
    36 = 1
    37 = 5
    38 = 0
    39 = 3
    40 = 0
    41 = 0
    42 = 1
    43 = 2
    44 = 0
    45 = 1
    46 = 4
    47 = 0
    48 = 2
    49 = 0
    50 = 1
    51 = 2
    52 = 2
    53 = 0
    54 = 2
    55 = 2
    56 = 0
    57 = 2
    58 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
